# daily auto push: 2026-01-05 02:38 UTC
# Append new rows 2875-2916 to Sheet1 (date/weekday/hour/ranking records
# for 2026/12/29 through 2027/01/05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(2875, "2026/12/29", "火", 13, 27),
    @(2876, "2026/12/29", "火", 16, 29),
    @(2877, "2026/12/29", "火", 19, 31),
    @(2878, "2026/12/29", "火", 23, 28),
    @(2879, "2026/12/30", "水", 2, 29),
    @(2880, "2026/12/30", "水", 5, 29),
    @(2881, "2026/12/30", "水", 8, 30),
    @(2882, "2026/12/30", "水", 13, 31),
    @(2883, "2026/12/30", "水", 16, 30),
    @(2884, "2026/12/30", "水", 22, 30),
    @(2885, "2026/12/31", "木", 2, 33),
    @(2886, "2026/12/31", "木", 6, 33),
    @(2887, "2026/12/31", "木", 9, 32),
    @(2888, "2026/12/31", "木", 12, 32),
    @(2889, "2026/12/31", "木", 14, 27),
    @(2890, "2026/12/31", "木", 22, 27),
    @(2891, "2027/01/01", "金", 2, 27),
    @(2892, "2027/01/01", "金", 5, 28),
    @(2893, "2027/01/01", "金", 13, 26),
    @(2894, "2027/01/01", "金", 16, 27),
    @(2895, "2027/01/01", "金", 19, 25),
    @(2896, "2027/01/02", "土", 1, 27),
    @(2897, "2027/01/02", "土", 5, 26),
    @(2898, "2027/01/02", "土", 8, 27),
    @(2899, "2027/01/02", "土", 13, 20),
    @(2900, "2027/01/02", "土", 16, 22),
    @(2901, "2027/01/02", "土", 19, 20),
    @(2902, "2027/01/02", "土", 22, 21),
    @(2903, "2027/01/03", "日", 1, 21),
    @(2904, "2027/01/03", "日", 4, 23),
    @(2905, "2027/01/03", "日", 7, 24),
    @(2906, "2027/01/03", "日", 13, 27),
    @(2907, "2027/01/03", "日", 16, 27),
    @(2908, "2027/01/03", "日", 19, 26),
    @(2909, "2027/01/03", "日", 22, 26),
    @(2910, "2027/01/04", "月", 2, 25),
    @(2911, "2027/01/04", "月", 5, 26),
    @(2912, "2027/01/04", "月", 7, 26),
    @(2913, "2027/01/04", "月", 13, 24),
    @(2914, "2027/01/04", "月", 22, 24),
    @(2915, "2027/01/05", "火", 1, 23),
    @(2916, "2027/01/05", "火", 7, 24)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A holds a date-looking string ("yyyy/mm/dd") that must stay a
    # literal text value (matches the rest of the column), not an Excel date
    # serial, so force text format before assigning, then restore the
    # default cell style so no extra formatting is left behind.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[1]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
